$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1697.8422
$ws.Range("J17").Value = 1697.8422
$ws.Range("L17").Value = 5093.5266
$ws.Range("N17").Value = -5429.5266
$ws.Range("H51").Value = 1765.3334
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 1765.3334
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 1765.3334
$ws.Range("N51").Value = -2733.3334
$ws.Range("M51").ClearContents()
$ws.Range("H125").Value = 5876
$ws.Range("I125").Value = 5732
$ws.Range("J125").Value = 5900
$ws.Range("K125").Value = 51588
$ws.Range("L125").Value = 53100
$ws.Range("M125").Value = -49128
$ws.Range("N125").Value = -58020
$ws.Range("H131").Value = 2959.9167
$ws.Range("I131").Value = 1541.4
$ws.Range("J131").Value = 10052.5
$ws.Range("K131").Value = 4624.200000000001
$ws.Range("L131").Value = 30157.5
$ws.Range("M131").Value = 415.7999999999993
$ws.Range("N131").Value = -40237.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2454.913
$ws.Range("I45").Value = 2223.9473
$ws.Range("J45").Value = 3552
$ws.Range("K45").Value = 2223.9473
$ws.Range("L45").Value = 3552
$ws.Range("M45").Value = -1846.9473
$ws.Range("N45").Value = -4306
$ws.Range("H122").Value = 868.1539
$ws.Range("I122").Value = 808.2222
$ws.Range("J122").Value = 1003
$ws.Range("K122").Value = 2424.6666
$ws.Range("L122").Value = 3009
$ws.Range("M122").Value = 25.33339999999998
$ws.Range("N122").Value = -7909

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 21900
$ws.Range("J6").Value = 21900
$ws.Range("L6").Value = 21900
$ws.Range("N6").Value = -22126
$ws.Range("H55").Value = 37750
$ws.Range("J55").Value = 37750
$ws.Range("L55").Value = 37750
$ws.Range("N55").Value = -38296
$ws.Range("H64").Value = 553.7143
$ws.Range("I64").Value = 1164.75
$ws.Range("J64").Value = 309.3
$ws.Range("K64").Value = 1164.75
$ws.Range("L64").Value = 309.3
$ws.Range("M64").Value = -939.75
$ws.Range("N64").Value = -759.3
$ws.Range("H67").Value = 553.7143
$ws.Range("I67").Value = 1164.75
$ws.Range("J67").Value = 309.3
$ws.Range("K67").Value = 1164.75
$ws.Range("L67").Value = 309.3
$ws.Range("M67").Value = -384.75
$ws.Range("N67").Value = -1869.3
$ws.Range("H114").Value = 37933.332
$ws.Range("J114").Value = 37933.332
$ws.Range("L114").Value = 37933.332
$ws.Range("N114").Value = -46611.332
$ws.Range("H123").Value = 47500
$ws.Range("J123").Value = 47500
$ws.Range("L123").Value = 47500
$ws.Range("N123").Value = -57300

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 3008.04
$ws.Range("J94").Value = 3924.1875
$ws.Range("L94").Value = 3924.1875
$ws.Range("N94").Value = -4826.1875
$ws.Range("H102").Value = 15000
$ws.Range("J102").Value = 15000
$ws.Range("L102").Value = 15000
$ws.Range("N102").Value = -19868
$ws.Range("H107").Value = 477045.75
$ws.Range("I107").Value = 667421.4399999999
$ws.Range("J107").Value = 1106.5
$ws.Range("K107").Value = 667421.4399999999
$ws.Range("L107").Value = 1106.5
$ws.Range("M107").Value = -665501.4399999999
$ws.Range("N107").Value = -4946.5
$ws.Range("H130").Value = 79700
$ws.Range("J130").Value = 79700
$ws.Range("L130").Value = 79700
$ws.Range("N130").Value = -89740
$ws.Range("H134").Value = 4119.8276
$ws.Range("I134").Value = 2302.8823
$ws.Range("J134").Value = 6693.8335
$ws.Range("K134").Value = 6908.646900000001
$ws.Range("L134").Value = 20081.5005
$ws.Range("M134").Value = -4373.646900000001
$ws.Range("N134").Value = -25151.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 571.4
$ws.Range("J5").Value = 679.3333
$ws.Range("L5").Value = 2037.9999
$ws.Range("N5").Value = -2261.9999
$ws.Range("H92").Value = 10000000
$ws.Range("I92").Value = 10000000
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 30000000
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -29998752
$ws.Range("N92").ClearContents()
$ws.Range("H122").Value = 7468.839
$ws.Range("I122").Value = 15526.357
$ws.Range("J122").Value = 833.2353000000001
$ws.Range("K122").Value = 139737.213
$ws.Range("L122").Value = 7499.117700000001
$ws.Range("M122").Value = -137287.213
$ws.Range("N122").Value = -12399.1177
$ws.Range("H131").Value = 1090.6177
$ws.Range("J131").Value = 1110.0312
$ws.Range("L131").Value = 3330.0936
$ws.Range("N131").Value = -13410.0936
$ws.Range("H134").Value = 4163.56
$ws.Range("I134").Value = 1826.2727
$ws.Range("K134").Value = 5478.8181
$ws.Range("M134").Value = -408.8181000000004
$ws.Range("H135").Value = 571.4
$ws.Range("J135").Value = 679.3333
$ws.Range("L135").Value = 6113.9997
$ws.Range("N135").Value = -11183.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 28850
$ws.Range("J39").Value = 28850
$ws.Range("L39").Value = 28850
$ws.Range("N39").Value = -29914
$ws.Range("H97").Value = 1432485.8
$ws.Range("I97").Value = 2503200
$ws.Range("J97").Value = 4866.6665
$ws.Range("K97").Value = 2503200
$ws.Range("L97").Value = 4866.6665
$ws.Range("M97").Value = -2502704
$ws.Range("N97").Value = -5858.6665
$ws.Range("H99").Value = 1979.6
$ws.Range("I99").Value = 1979.6
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1979.6
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 266.4000000000001
$ws.Range("N99").ClearContents()
$ws.Range("H102").Value = 1520.6
$ws.Range("I102").Value = 1518.909
$ws.Range("K102").Value = 1518.909
$ws.Range("M102").Value = 103.0909999999999
$ws.Range("H132").Value = 4584.364
$ws.Range("I132").Value = 4912.6284
$ws.Range("J132").Value = 3307.7778
$ws.Range("K132").Value = 14737.8852
$ws.Range("L132").Value = 9923.3334
$ws.Range("M132").Value = -12207.8852
$ws.Range("N132").Value = -14983.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 5344.3335
$ws.Range("I43").Value = 2033
$ws.Range("J43").Value = 7000
$ws.Range("K43").Value = 2033
$ws.Range("L43").Value = 7000
$ws.Range("M43").Value = -1884
$ws.Range("N43").Value = -7298
$ws.Range("H49").Value = 7500
$ws.Range("J49").Value = 7500
$ws.Range("L49").Value = 7500
$ws.Range("N49").Value = -7960
$ws.Range("H122").Value = 3231.7334
$ws.Range("I122").Value = 2593.7144
$ws.Range("J122").Value = 3790
$ws.Range("K122").Value = 7781.1432
$ws.Range("L122").Value = 11370
$ws.Range("M122").Value = -5331.1432
$ws.Range("N122").Value = -16270

